# Generate Report for Handback
#
# This localization-status report previously showed every locale as
# "Ready for handoff" with no handback info recorded yet. This pass marks
# the status as handed back, records the target/handback files that were
# produced, and stamps the handback datetime - for zh-cn (still pending
# translation, so it is simply back "in sync with en-US") and for de-de
# (which now has a real handback package).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$u871 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ba14010197aeac15c0a89008f6e37e09a65b7cb8/e2e/871c2a8f-e428-4c32-a3d7-4f079e4772ac.md"
$ua278 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ba14010197aeac15c0a89008f6e37e09a65b7cb8/e2e/a278e830-8c0f-430d-b262-e3153c0360f6.md"

$d871 = "871c2a8f-e428-4c32-a3d7-4f079e4772ac.md"
$da278 = "a278e830-8c0f-430d-b262-e3153c0360f6.md"

# ---------------------------------------------------------------------
# Overview sheet: Status columns (E, F) for both rows flip to the new
# "handed back" status text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Widen the (now longer) status columns to fit the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet: status + target/handback file columns + handback datetime.
# zh-cn only got re-synced with en-US, so its handback datetime is
# 2016-08-17 10:58:29.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsZh.Range("J2").Value = "871c2a8f-e428-4c32-a3d7-4f079e4772ac.f50b0bc230aa48e54d32b4d2ccdadeed5a020b4a.zh-cn.xlf"
$wsZh.Range("J3").Value = "a278e830-8c0f-430d-b262-e3153c0360f6.8788f3273d79f0ab53795567bf8e8702f5add239.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-17 10:58:29"
$wsZh.Range("K3").Value = "2016-08-17 10:58:29"

# Rebuild the hyperlinks in display order (A2, I2, A3, I3) so the new
# "Latest Target File" links (I2/I3) land next to the existing "Source
# File Name" links (A2/A3) pointing at the same two files.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $u871, "", "", $d871)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $u871, "", "", $d871)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $ua278, "", "", $da278)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $ua278, "", "", $da278)

# ---------------------------------------------------------------------
# de-de sheet: status + target/handback file columns + handback datetime.
# de-de has a fresh handback package, stamped 2016-08-17 10:58:36.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDe.Range("J2").Value = "871c2a8f-e428-4c32-a3d7-4f079e4772ac.f50b0bc230aa48e54d32b4d2ccdadeed5a020b4a.de-de.xlf"
$wsDe.Range("J3").Value = "a278e830-8c0f-430d-b262-e3153c0360f6.8788f3273d79f0ab53795567bf8e8702f5add239.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-17 10:58:36"
$wsDe.Range("K3").Value = "2016-08-17 10:58:36"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $u871, "", "", $d871)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $u871, "", "", $d871)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $ua278, "", "", $da278)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $ua278, "", "", $da278)
